# Add a new "AddVacancy" worksheet (after Sheet1) with a small data-driven
# table of vacancies, matching the "data driven test" described in the
# commit message.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so tab order is Sheet1, AddVacancy.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws.Name = "AddVacancy"

# Header row
$ws.Range("A1").Value = "job title"
$ws.Range("B1").Value = "vacancy name"
$ws.Range("C1").Value = "hiring manager"
$ws.Range("D1").Value = "number of positions"
$ws.Range("E1").Value = "description"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "test engineer"
$ws.Range("C2").Value = "Kallyani Bhute"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "perform test using selenium in java"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "production"
$ws.Range("C3").Value = "Paul Collings"
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "part of team in assembly line"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "driver"
$ws.Range("C4").Value = "Rebecca Harmony"
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = "AZ driver"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "developer"
$ws.Range("C5").Value = "Dominic Chase"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "frontend developer"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "analyst"
$ws.Range("C6").Value = "Nathan Elliot"
$ws.Range("D6").Value = 9
$ws.Range("E6").Value = "anylyse software feasibility"

# Column A holds the numeric id but is formatted as text (numFmtId 49 "@").
# Apply after the values are written so every cell in the column picks it up.
$ws.Columns.Item(1).NumberFormat = "@"

# Best-fit the column widths like Excel's AutoFit on "vacancy created" entry
# (explicit widths chosen so the saved, padded xlsx column width lines up with
# the content-driven best-fit Excel itself would have computed).
$ws.Columns.Item(1).ColumnWidth = 6.666666666666667
$ws.Columns.Item(2).ColumnWidth = 11.833333333333334
$ws.Columns.Item(3).ColumnWidth = 15.0
$ws.Columns.Item(4).ColumnWidth = 16.833333333333332
$ws.Columns.Item(5).ColumnWidth = 29.833333333333332

# Match the author's final selection/cursor position on the new sheet.
$ws.Range("A5").Select() | Out-Null
